$wb = $excel.ActiveWorkbook

# "Restricciones_del_follower" sheet (3rd sheet) - columns E (J_0_LP_v) and F (Gamma_value)
$wsFollower = $wb.Worksheets.Item(3)
$wsFollower.Cells.Item(2, 5).Value = "'1.7000000000000002"
$wsFollower.Cells.Item(2, 6).Value = "'2.4"
$wsFollower.Cells.Item(3, 5).Value = "'0.3"
$wsFollower.Cells.Item(3, 6).Value = "'3.7"

# "Vector_BF" sheet (6th sheet) - cell A3
$wsVectorBF = $wb.Worksheets.Item(6)
$wsVectorBF.Cells.Item(3, 1).Value = "'-68.85"
